# Disaggregation of commodity Copper
#
# 1) Rename the "Copper ores and concentrates" commodity to "Copper"
#    (row 7, column C, on every year sheet).
# 2) The historic values that used to live under "Copper ores and
#    concentrates" are redistributed across the three end-use columns
#    (D = Photovoltaic plants, E = Offshore wind plants, F = Onshore wind
#    plants) for rows 5 (Neodymium), 7 (Copper) and 8 (Raw silicon): each
#    row's three values are cycled one column to the right, wrapping
#    around (new D = old F, new E = old D, new F = old E).

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Rename the commodity label in column C, row 7.
    $ws.Cells.Item(7, 3).Value2 = "Copper"

    # Rotate D/E/F values for rows 5, 7 and 8.
    foreach ($r in 5, 7, 8) {
        $dVal = $ws.Cells.Item($r, 4).Value2
        $eVal = $ws.Cells.Item($r, 5).Value2
        $fVal = $ws.Cells.Item($r, 6).Value2

        $ws.Cells.Item($r, 4).Value2 = $fVal
        $ws.Cells.Item($r, 5).Value2 = $dVal
        $ws.Cells.Item($r, 6).Value2 = $eVal
    }
}
